$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value map, derived from the authoritative diff of the
# refreshed crypto price/volume feed. Every value must be written back as
# literal text (matching the original inlineStr cells), so we temporarily
# force a Text number format before assigning, then restore the default
# 'Normal' style so no stray formatting is introduced.
$updates = @(
    @{ Cell = 'D2'; Value = '67.118.25' }
    @{ Cell = 'E2'; Value = '  +4.64%  ' }
    @{ Cell = 'D3'; Value = '3.467.96' }
    @{ Cell = 'E3'; Value = '  +4.34%  ' }
    @{ Cell = 'E4'; Value = '  +0.08%  ' }
    @{ Cell = 'D5'; Value = '585.94' }
    @{ Cell = 'E5'; Value = '  +6.36%  ' }
    @{ Cell = 'D6'; Value = '188.44' }
    @{ Cell = 'E6'; Value = '  +9.03%  ' }
    @{ Cell = 'D7'; Value = '0.632' }
    @{ Cell = 'E7'; Value = '  +1.09%  ' }
    @{ Cell = 'D8'; Value = '3.462.68' }
    @{ Cell = 'E8'; Value = '  +4.52%  ' }
    @{ Cell = 'D10'; Value = '0.172' }
    @{ Cell = 'E10'; Value = '  +0.23%  ' }
    @{ Cell = 'D11'; Value = '0.649' }
    @{ Cell = 'E11'; Value = '  +2.24%  ' }
    @{ Cell = 'D12'; Value = '56.87' }
    @{ Cell = 'E12'; Value = '  +7.02%  ' }
    @{ Cell = 'E13'; Value = '  +0.42%  ' }
    @{ Cell = 'D14'; Value = '9.43' }
    @{ Cell = 'E14'; Value = '  +4.33%  ' }
    @{ Cell = 'D15'; Value = '4.027.50' }
    @{ Cell = 'E15'; Value = '  +4.50%  ' }
    @{ Cell = 'D16'; Value = '18.76' }
    @{ Cell = 'E16'; Value = '  +4.00%  ' }
    @{ Cell = 'D17'; Value = '3.471.95' }
    @{ Cell = 'E17'; Value = '  +4.12%  ' }
    @{ Cell = 'D18'; Value = '67.149.17' }
    @{ Cell = 'E18'; Value = '  +4.92%  ' }
    @{ Cell = 'D19'; Value = '12.16' }
    @{ Cell = 'E19'; Value = '  +4.13%  ' }
    @{ Cell = 'D20'; Value = '0.118' }
    @{ Cell = 'E20'; Value = '  -1.44%  ' }
    @{ Cell = 'E21'; Value = '  +4.07%  ' }
    @{ Cell = 'D22'; Value = '485.99' }
    @{ Cell = 'E22'; Value = '  +8.49%  ' }
    @{ Cell = 'D23'; Value = '5.37' }
    @{ Cell = 'E23'; Value = '  +7.95%  ' }
    @{ Cell = 'D24'; Value = '16.83' }
    @{ Cell = 'E24'; Value = '  +21.36%  ' }
    @{ Cell = 'D25'; Value = '4.47' }
    @{ Cell = 'E25'; Value = '  +10.83%  ' }
    @{ Cell = 'D26'; Value = '89.65' }
    @{ Cell = 'E26'; Value = '  +3.18%  ' }
    @{ Cell = 'E27'; Value = '  +3.12%  ' }
    @{ Cell = 'D28'; Value = '10.95' }
    @{ Cell = 'E28'; Value = '  +3.69%  ' }
    @{ Cell = 'D29'; Value = '9.11' }
    @{ Cell = 'E29'; Value = '  +6.47%  ' }
    @{ Cell = 'D30'; Value = '31.38' }
    @{ Cell = 'E30'; Value = '  +2.08%  ' }
    @{ Cell = 'D31'; Value = '7.18' }
    @{ Cell = 'E31'; Value = '  +10.76%  ' }
    @{ Cell = 'B32'; Value = 'Cosmos' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' }
    @{ Cell = 'D32'; Value = '11.75' }
    @{ Cell = 'E32'; Value = '  +3.63%  ' }
    @{ Cell = 'B33'; Value = 'Bittensor' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = 'D33'; Value = '596.76' }
    @{ Cell = 'E33'; Value = '  +4.68%  ' }
    @{ Cell = 'D34'; Value = '64.10' }
    @{ Cell = 'E34'; Value = '  +2.45%  ' }
    @{ Cell = 'E35'; Value = '  +5.27%  ' }
    @{ Cell = 'E36'; Value = '  +6.71%  ' }
    @{ Cell = 'E37'; Value = '  -0.03%  ' }
    @{ Cell = 'D38'; Value = '36.68' }
    @{ Cell = 'E38'; Value = '  +4.47%  ' }
    @{ Cell = 'D39'; Value = '3.54' }
    @{ Cell = 'E39'; Value = '  +0.13%  ' }
    @{ Cell = 'E40'; Value = '  +5.42%  ' }
    @{ Cell = 'D41'; Value = '0.0₃0758' }
    @{ Cell = 'E41'; Value = '  +4.55%  ' }
    @{ Cell = 'D42'; Value = '3.236.41' }
    @{ Cell = 'E42'; Value = '  +6.00%  ' }
    @{ Cell = 'D43'; Value = '2.91' }
    @{ Cell = 'E43'; Value = '  +7.03%  ' }
    @{ Cell = 'E44'; Value = '  +4.50%  ' }
    @{ Cell = 'B45'; Value = 'dogwifhat' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' }
    @{ Cell = 'D45'; Value = '2.79' }
    @{ Cell = 'E45'; Value = '  +24.79%  ' }
    @{ Cell = 'B46'; Value = 'ApeXProtocol' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex' }
    @{ Cell = 'D46'; Value = '3.27' }
    @{ Cell = 'E46'; Value = '  +3.35%  ' }
    @{ Cell = 'B47'; Value = 'Fetch.AI' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D47'; Value = '2.51' }
    @{ Cell = 'E47'; Value = '  +3.23%  ' }
    @{ Cell = 'D48'; Value = '0.135' }
    @{ Cell = 'E48'; Value = '  +1.65%  ' }
    @{ Cell = 'D49'; Value = '3.27' }
    @{ Cell = 'E49'; Value = '  +13.68%  ' }
    @{ Cell = 'D50'; Value = '8.74' }
    @{ Cell = 'E50'; Value = '  +7.11%  ' }
    @{ Cell = 'D51'; Value = '1.00' }
    @{ Cell = 'E51'; Value = '  +0.13%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

